$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Generate Report for Handback
#
# For each per-language sheet (zh-cn, de-de):
#   - Status column (C) flips from "Ready for handoff" to
#     "Handed back: in sync with en-US" for both data rows.
#   - Two new columns are populated for both data rows:
#       F = "Latest Target File"    (same file identity as column A, the .md)
#       G = "Latest Handback File"  (same file identity as column D, the .xlf)
#     each carrying a hyperlink, mirroring the A/D hyperlink pattern.
#   - Latest Handback DateTime (H) is stamped with the handback completion
#     time (distinct per sheet, since zh-cn finished before de-de).
# ---------------------------------------------------------------------------

$statusText = "Handed back: in sync with en-US"

$sheets = @(
    @{
        Name = "zh-cn"
        Row2Base = "3c57661d-4743-4840-9c7c-46745d6a599a"
        Row2Hash = "a46db04ae91046f46f49fe1bd38db7d67d0decbe"
        Row3Base = "d1a30974-23e3-45dc-acce-fb657c1ffbf5"
        Row3Hash = "ec3097626fa649c88075ba5cd1f172b135ddd8e4"
        Lang = "zh-cn"
        HandbackDateTime = "2016-03-19 04:03:28"
        MdRoot = "https://github.com/OpenLocalizationTest/oltest/blob/48c463dfd0aff8820597b6fed29253fd55971103/e2e"
        XlfRoot = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9f3be9e2ea86cb360d204200594b521ee80b947a/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht"
    },
    @{
        Name = "de-de"
        Row2Base = "3c57661d-4743-4840-9c7c-46745d6a599a"
        Row2Hash = "a46db04ae91046f46f49fe1bd38db7d67d0decbe"
        Row3Base = "d1a30974-23e3-45dc-acce-fb657c1ffbf5"
        Row3Hash = "ec3097626fa649c88075ba5cd1f172b135ddd8e4"
        Lang = "de-de"
        HandbackDateTime = "2016-03-19 04:03:41"
        MdRoot = "https://github.com/OpenLocalizationTest/oltest/blob/48c463dfd0aff8820597b6fed29253fd55971103/e2e"
        XlfRoot = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/7f01d7d32798d2116750edb515c5f0a481bafea5/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht"
    }
)

foreach ($s in $sheets) {
    $ws = $wb.Worksheets.Item($s.Name)

    # --- Status column (C2:C3) -> handed back ---------------------------
    $ws.Range("C2").Value = $statusText
    $ws.Range("C3").Value = $statusText

    # --- row 2 (file 3c57661d-...) ---------------------------------------
    $mdDisplay2  = "$($s.Row2Base).md"
    $xlfDisplay2 = "$($s.Row2Base).$($s.Row2Hash).$($s.Lang).xlf"

    $ws.Range("F2").Value = $mdDisplay2
    $ws.Range("F2").Font.Underline = $true
    $ws.Range("F2").Font.Color = 13272813
    $ws.Hyperlinks.Add($ws.Range("F2"), "$($s.MdRoot)/$mdDisplay2", "", "", $mdDisplay2) | Out-Null

    $ws.Range("G2").Value = $xlfDisplay2
    $ws.Range("G2").Font.Underline = $true
    $ws.Range("G2").Font.Color = 13272813
    $ws.Hyperlinks.Add($ws.Range("G2"), "$($s.XlfRoot)/$xlfDisplay2", "", "", $xlfDisplay2) | Out-Null

    # --- row 3 (file d1a30974-...) -----------------------------------------
    $mdDisplay3  = "$($s.Row3Base).md"
    $xlfDisplay3 = "$($s.Row3Base).$($s.Row3Hash).$($s.Lang).xlf"

    $ws.Range("F3").Value = $mdDisplay3
    $ws.Range("F3").Font.Underline = $true
    $ws.Range("F3").Font.Color = 13272813
    $ws.Hyperlinks.Add($ws.Range("F3"), "$($s.MdRoot)/$mdDisplay3", "", "", $mdDisplay3) | Out-Null

    $ws.Range("G3").Value = $xlfDisplay3
    $ws.Range("G3").Font.Underline = $true
    $ws.Range("G3").Font.Color = 13272813
    $ws.Hyperlinks.Add($ws.Range("G3"), "$($s.XlfRoot)/$xlfDisplay3", "", "", $xlfDisplay3) | Out-Null

    # --- Latest Handback DateTime (H2:H3) -----------------------------------
    $ws.Range("H2").Value = $s.HandbackDateTime
    $ws.Range("H3").Value = $s.HandbackDateTime
}
